# Load display.pptx - "Adding the load display subsystem diagram"
#
# 1) The "datetimeFigureOut" Date placeholder on the Slide Master and on
#    every Slide Layout had its cached text bumped from 3/13/2014 to
#    3/19/2014 (the field's cached/"last computed" text, as happens when
#    the deck is re-saved on a later date).
# 2) The "TextBox 32" caption shape on slide 1 shifted left (x offset
#    1362021 -> 1325534 EMU; y/width/height unchanged).

$p = $ppt.ActivePresentation

$oldDateText = "3/13/2014"
$newDateText = "3/19/2014"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)

        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            # msoPlaceholder
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDateText) {
                $tr.Text = $newDateText
            }
        }
    }
}

# Slide Master's own Date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout's Date placeholder.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Reposition the "TextBox 32" caption on slide 1.
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.Name -eq "TextBox 32") {
        # 1325534 EMU (was 1362021 EMU); 914400 EMU per inch, 12700 EMU per point.
        $shp.Left = 1325534 / 12700
    }
}
